$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new "Building Description" data in column G (donor biography text)
# and fill in "Null" placeholders in column H (GPS Coordinates) where no coordinates exist yet.

$ws.Range('G2').Value = 'Robert C. Pew was the distinguished founder and President of Steelcase from 1966 to 1974.  He was a prominent member of the citizen’s group that established the GVSU charter.  He was also instrumental in chairing the committee that raised over nine million dollars to build and equip the Eberhard Center.   He was a member of the board of trustees from 1985 to 1988 as well as donating 10 acres of real estate that were pivotal in building the Grand Rapids campus the bears his name.  He was a member of the board of trustees.'
$ws.Range('G2').VerticalAlignment = -4108
$ws.Range('H2').Value = 'Null'

$ws.Range('G3').Value = 'Thomas E. Beckering was the third generation owner of Pioneer Construction.  He was influential in propelling Pioneer Construction to become one of the largest full-service construction companies in the Midwest.'
$ws.Range('G3').VerticalAlignment = -4108
$ws.Range('H3').Value = 'Null'

$ws.Range('G4').Value = 'Peter C. Cook created Mazda Great Lakes, the first and only independent direct distributorship in the U.S. until Mazda bought it back.  He is a renowned member of the prestigious Grand Valley Hall of Fame.  He also enjoyed playing the trumpet, as well as he had a love for Dixie bands.  He believed strongly in giving back to his community stating that, “There are two necessary things in giving: one is to be willing, the other is to be able.  Able is secondary because many people who are not willing to do become able but still never give.  People who have little yet give what they can become the big givers because they are able.  That’s because they had always been willing.”'
$ws.Range('G4').VerticalAlignment = -4108

$ws.Range('G5').Value = 'Richard M. DeVos was the president and chairman of the Grand valley University Foundation.  He is the co-founder of Amway Corporation.  He is both an entrepreneur and a philanthropist.  He has a love for sports, owning the Orlando Magic, and NBA franchise in Florida, naming the arena for the company he started.  He was instrumental in leading the revitalization of downtown Grand Rapids when he opened the Amway Grand plaza Hotel.  He is one of Grand Valleys’ most charitable donors leading nearly every campaign since the early 1990’s.  '
$ws.Range('G5').VerticalAlignment = -4108

$ws.Range('G6').Value = 'L.V Eberhard is the founder of Eberhard Foods, a company who spent more the 70 years in the grocery business.  He donated to Grand Valley in 1988 for the expansion of the downtown campus as well as financing addition scholarship for students in need.  He was a member of the Grand Valley University Foundation Advisory Cabinet form 1989 to 1991.  He is a noted charitable member of our community, receiving the Enrichment Award in 1988.'
$ws.Range('G6').VerticalAlignment = -4108

$ws.Range('G7').Value = 'L. William “Bill” Seidman is the Founder of Grand Valle State University, leading the efforts to charter a public university in West Michigan.  He served as the board of control chairman for two separate sessions, both between 1960 and 1974 as well as between 1977 and 1983.  He is a member of the Grand Valley State University Hall of Fame.  And he even served as the chairman of the FDIC between 1985 and 1991.'
$ws.Range('G7').VerticalAlignment = -4108
$ws.Range('H7').Value = 'Null'

$ws.Range('G8').Value = 'Fred M. Keller is the individual who patented the single-action flush lever.   As a business man he bout Paragon Tool & Die Company for a single dollar, he would later become a member of the Greater Grand Rapids Business Hall of Fame.   Although being recognized for his charitable contributions to our community by receiving the Enrichment Award in 2000, he is not the only charitable one in his family.  Both Fred and his wife, Bernie believed in giving back to their community and were both recognized as two of the community’s largest philanthropists.'
$ws.Range('G8').VerticalAlignment = -4108

$ws.Range('G9').Value = 'John C. Kennedy III is currently a member of the Grand Valley State University Board of Trustees, starting his term in 2010.  He founded Autocam in 1988 when he was only thirty.  He is a member of the West Michigan Business Hall of Fame.  He remains highly involve in West Michigan, including Grand Rapids Catholic Secondary Schools as well as Saint Mary’s Health Care.  Both his wife Nancy and he have received the Chaille Community Philanthropy Award from the Grand Rapids Community Foundation in 2012.  He is also a recipient of the Grand Valley University Foundation Leadership Award in 2002.'
$ws.Range('G9').VerticalAlignment = -4108

$ws.Range('G10').Value = 'Peter F. Secchia served as the United States ambassador to Italy between 1989 and 1993, and was declared Michigan’s Entrepreneur of the Year in 1994.  He received an honorary Doctor of Laws Degree from GVSU in 1991.  He also donated the bronze sculpture of Chief Noonday that is located near the blue bridge.'
$ws.Range('G10').VerticalAlignment = -4108

$ws.Range('G11').Value = 'Richard M. DeVos was the president and chairman of the Grand valley University Foundation.  He is the co-founder of Amway Corporation.  He is both an entrepreneur and a philanthropist.  He has a love for sports, owning the Orlando Magic, and NBA franchise in Florida, naming the arena for the company he started.  He was instrumental in leading the revitalization of downtown Grand Rapids when he opened the Amway Grand plaza Hotel.  He is one of Grand Valleys’ most charitable donors leading nearly every campaign since the early 1990’s.  '
$ws.Range('G11').VerticalAlignment = -4108

$ws.Range('G12').Value = 'Frederik Meijer worked in his father’s grocery store, and together they created the first “supercenter” concept with Meijer thrifty acres in 1962.  He created the Fredrik Meijer gardens and sculpture park in in Grand Rapids.  At Grand Valley he also helped establish the Fredrik Meijer Honors College as well as the Meijer campus in Holland.  He served as the Grand Valley University Foundation Director between 1981 and 2011. He had a love for the arts, education, health and the many parks of West Michigan.'
$ws.Range('G12').VerticalAlignment = -4108
$ws.Range('H12').Value = 'Null'

$ws.Range('G13').Value = 'Steelcase was founded in Grand Rapids in 1912 and has become one of the world’s largest furniture manufacturers.  Steelcase donated a substantial portion of land to Grand Valley that became part of the Robert C. Pew campus.  As a company they have always believed in supporting their community through their beliefs and values stating that, “being a good neighbor is just an important part of being a good company.”  Like there founder, Robert C. Pew, as a company they have always been supportive of the arts, the environment, human services and economic development throughout West Michigan.'

$ws.Range('G14').Value = 'Alexander Calder is the American sculptor credited with inventing the mobile.  He also created prolific paintings, lithographs, toys, tapestry and jewelry.  One of his most famous sculptures is the, “La Grand Vitesse” which still stands in downtown Grand Rapids.  La Grand Vitesse was the first public art in the United States that was funded with federal money, being financed by the National Endowment for the Arts.'
$ws.Range('G14').VerticalAlignment = -4108

$ws.Range('G15').Value = 'Arthur C. Hills was part of the Original faculty at Grand Valley State University when the school opened in 1958.  In 1963 he was charged with creating an arts curriculum and organize arts activities for all students.  In the spring of 1964 he produced the fine arts on campus series for students and community.  He is renowned for creating the Grand Valley alma mater, “Hail to Thee Grand Valley.”  He remained an active member of the Grand Valley community till the day he retired in 1988.'
$ws.Range('G15').VerticalAlignment = -4108

$ws.Range('G16').Value = 'Icie Macy Hoobler served as a member of the Grand Valley State College Board of Trustees between 1960 and 1968.  She was a scientist, an educator, a leader, an author and much more.  She was the first woman president of the American Chemical Society.  Her works included heading the Nutrition Research Project for the Merrill-Palmer School in Detroit, where she studied the effects of various foods on growth and well-being of children.  She is an honorary recipient of a Grand Valley doctorate as well as Humane Letters.'
$ws.Range('G16').VerticalAlignment = -4108

$ws.Range('G17').Value = 'James M. Copeland served as a member of the Grand Valley State College Board of Trustees between 1960 and 1968 where he was instrumental in starting the Grand Valley State College Library.  As a businessman he was the President of Security National Bank in Manistee as well as the President of the Wyoming State Bank in Grand Rapids.'
$ws.Range('G17').VerticalAlignment = -4108

$ws.Range('G18').Value = 'Edward J. Frey served on the Grand Valley State College founding committee, and was instrumental in selecting the first president, James H. Zumberge.  He is the heir to the Union Bank and Trust Company.  Frey and his wife Frances we both civic leaders in West Michigan, establishing the Frey Foundation in 1974 to help fund the arts, education, economic development, and to help improve the lives of women and children.  He also helped establish Grand Rapids Chapter of Junior Achievement.'
$ws.Range('G18').VerticalAlignment = -4108

$ws.Range('G19').Value = 'Paul A. Johnson served on the Grand Valley State College founding committee, as well as serving on the Grand Valley University Foundation between 1981 and 2003.  He was a member of the Grand Valley State University Board of Trustees between 1968 and 1994.  In his business life he acted as the chairman of JSJ Corporation, and also the Kent County Director of Road Maintenance between 1978 and 1986.  HE would also receive an Honorary Doctorate in Humanities in 1994.'
$ws.Range('G19').VerticalAlignment = -4108

$ws.Range('G20').Value = 'Brian Kelley served as GVSU’s Head Football Coach from 1991-2003 and in 2009 was inducted into the Laker Hall of Fame.  While acting as head coach he led GVSU to National Championships in 2002 & 2003.  Paqui, Brian’s wife, met Brian while she was working at GVSU in Financial Aid.  And together they started the Kelly Cares Foundation in 2007.'
$ws.Range('G20').VerticalAlignment = -4108

$ws.Range('G21').Value = 'Russel H. Kirkhof donated a million dollars to Grand Valley State College in 1978 to help maintain the landscaping throughout campus.  The Kirkhof Center and Kirkhof College of Nursing were both named in honor of this gift, which at the time was the largest gift ever received by the university.'
$ws.Range('G21').VerticalAlignment = -4108

$ws.Range('G22').Value = 'William A. Kirkpatrick served as a member of the Grand Valley State College Board of Trustees between 1960 and 1978.  In his business life he acted as President, Treasurer, and Board Member, Kalamazoo Paper Box Company.  As well as serving as President of Goodwill Industries—Southwestern Michigan.'
$ws.Range('G22').VerticalAlignment = -4108

$ws.Range('G23').Value = 'Robert Kleiner served on the Grand Valley State University Board of Trustees between 1983 and 1987.  He was a practicing attorney as well as a supporter of Democratic Party, acting as a Democratic National Convention Delegate in 1964, 1968 and 1976.  With a love for art he was also a member of the Grand Rapids Art Museum Board.'

$ws.Range('G24').Value = 'Grace Olsen Kistler was a member of the Grand Valley State College Board of Control between 1960 and 1967.  She was a devoted volunteer and community leader, serving as a Board Member of the Ottawa County Red Cross.  She was the president of the State Federation of Women’s Clubs in Michigan.  And she even served as the secretary of the North Ottawa Community Hospital Authority.'
$ws.Range('G24').VerticalAlignment = -4108

$ws.Range('G25').Value = 'The Loutit Foundation was formed in 1957 by William R. Loutit.  The Loutit family has ties to Grand Haven and Spring Lake that date back to 1856 when Captain William R. Loutit settled in the area and grew a lumber and shipping business.  They were huge supporters of Grand Valley State College’s construction of a new science building—GVSU was first and only public institution to receive funding from the foundation.'
$ws.Range('G25').VerticalAlignment = -4108

$ws.Range('G26').Value = 'Arend D. Lubbers served as the Grand Valley State University President between 1969 and 2001, becoming president of GVSU at only 37 years old.  He was very influential in establishing the Seidman College of Business, Kirkhof College of Nursing, Schools of Education, Social Work, Health Sciences, Engineering, and Communication during his tenure.   Famous nation-wide LIFE Magazine named him to list of top 100 young men under 40 in the U.S.  He is also a Laker Hall of Fame Inductee.  '
$ws.Range('G26').VerticalAlignment = -4108

$ws.Range('G27').Value = 'Mark A. Murray was the third President of Grand Valley State University serving between 2001 and 2006.  Participating in politics in his career he also served as Treasurer of Michigan under Governor John Engler.  As an educator he acted as Michigan Education Assessment Program (MEAP) testing and developed the Michigan Merit Award.  In 2006 he would move on to become President of Meijer, Inc.'
$ws.Range('G27').VerticalAlignment = -4108

$ws.Range('G28').Value = 'Glenn A. Niemeyer was one of first faculty members appointed at Grand Valley in 1963, serving as the Dean of the College of Arts & Sciences between 1970 and 1976, as well as the Vice President for Academic Affairs between 1976 and 1980, and even as Provost between 1980 and 2001.  He chaired the North Central Association Higher Learning Commission, as well as acting as Vice President of North Central Association of Colleges and Schools.'
$ws.Range('G28').VerticalAlignment = -4108

$ws.Range('G29').Value = 'Glenn A. Niemeyer was one of first faculty members appointed at Grand Valley in 1963, serving as the Dean of the College of Arts & Sciences between 1970 and 1976, as well as the Vice President for Academic Affairs between 1976 and 1980, and even as Provost between 1980 and 2001.  He chaired the North Central Association Higher Learning Commission, as well as acting as Vice President of North Central Association of Colleges and Schools.'
$ws.Range('G29').VerticalAlignment = -4108

$ws.Range('G30').Value = 'Glenn A. Niemeyer was one of first faculty members appointed at Grand Valley in 1963, serving as the Dean of the College of Arts & Sciences between 1970 and 1976, as well as the Vice President for Academic Affairs between 1976 and 1980, and even as Provost between 1980 and 2001.  He chaired the North Central Association Higher Learning Commission, as well as acting as Vice President of North Central Association of Colleges and Schools.'
$ws.Range('G30').VerticalAlignment = -4108

$ws.Range('G31').Value = 'Arnold C. Ott served on Grand Valley’s Board of Trustees for 29 years, which is the longest serving board member at a public university in Michigan.  He founded the Ott Chemical Company in Muskegon in 1956.  He was a prolific inventor authoring 50+ domestic and foreign patents, even acting as a Science Consultant for the Amway Corporation.  HE would eventually establish the Arnold C. Ott Lectureship in Chemistry in 2001 at GVSU.'
$ws.Range('G31').VerticalAlignment = -4108

$ws.Range('G32').Value = 'Seymour Padnos is the Chairman and CEO, Louis Padnos Iron & Metal Company, 1955-Present.  He was inducted into GVSU Hall of Fame in 1988 as well as the West MI Business Hall of Fame.  He served as trustee to the Gerald R. Ford Foundation.  Both Seymore and Esther, Seymores Wife, were huge supporters of projects for the State of Michigan and the Holland area.'
$ws.Range('G32').VerticalAlignment = -4108

$ws.Range('G33').Value = 'Robert C. Pew was the distinguished founder and President of Steelcase from 1966 to 1974.  He was a prominent member of the citizen’s group that established the GVSU charter.  He was also instrumental in chairing the committee that raised over nine million dollars to build and equip the Eberhard Center.   He was a member of the board of trustees from 1985 to 1988 as well as donating 10 acres of real estate that were pivotal in building the Grand Rapids campus the bears his name.  He was a member of the board of trustees.'
$ws.Range('G33').VerticalAlignment = -4108

$ws.Range('G34').Value = 'William F. Pickard was born in La Grange, Georgia and moved to Flint, MI in 1958.  He became active in NAACP and was well-known leader in the black community.  He was successful on many fronts both as an entrepreneur and college teacher.  He moved to Detroit to oversee McDonald’s franchise in 1972, eventually acquiring six restaurants.  He joined the GVSU Board of Control and served between 1976 and 1988.  In 1972 he would receive an Honorary doctorate from GVSU.'
$ws.Range('G34').VerticalAlignment = -4108

$ws.Range('G35').Value = 'Kenneth W. Robinson was born in Grand Rapids, attended Central High School.  He was the Regional Director of UAW between 1955 and 1974.  He was active both in politics and civil rights.  He was even a Delegate at several Democratic National Conventions.  He was appointed to Grand Valley State College’s first Board of Control in 1960 and reappointed in 1968, 1972, 1973.'
$ws.Range('G35').VerticalAlignment = -4108

$ws.Range('G36').Value = 'L. William “Bill” Seidman is known as Grand Valley State University’s founding father.  He led the initiative to charter Grand Valley State College in 1960.  He was a partner of Seidman & Seidman – accounting firm founded by his father, Frank Edward Seidman (namesake of Seidman College of Business).  He served as Special Assistant on Financial Affairs to Governor George Romney between 1963 and 1966.  He also served as an economic advisor to President Gerald R. Ford between 1974 and 1977.  He chaired the Federal Deposit Insurance Corporation (FDIC) under President Ronald Reagan and George H.W. Bush.  He is credited with being very influential in leading  the US through banking/financial crisis.'
$ws.Range('G36').VerticalAlignment = -4108

$ws.Range('G37').Value = 'Dale Stafford was an accomplished journalist and newspaper publisher working for the Pontiac Press, Associated Press and the Detroit Free Press.  He was a correspondent for the Lansing State Journal and sports writer for Lansing Capital News while in high school.  He was the Editor and Publisher, Greenville Daily News, as well as the president of the Greenville Chamber of Commerce.   He was the first publicity director for athletics at Michigan State College.   He was appointed to Grand Valley’s Board of Control and served between 1960 and 1968.  '
$ws.Range('G37').VerticalAlignment = -4108

$ws.Range('G38').Value = 'Maxine M. Swanson was a Nurse and administrator at Gratiot Community Hospital.  She was very active with GOP in her county – began as precinct delegate, then Gratiot County Executive Committee, then chairman of President Ford campaign.  She was also a member the League of Women Voters, Coalition to Women’s Assembly, Michigan Republican Women’s Task Force, United Way, Alma Symphony, Sierra Club, Michigan Society for Mental Health.  She also served as a member of the Grand Valley State University Board of Trustees between 1981 and 1990.'
$ws.Range('G38').VerticalAlignment = -4108

$ws.Range('G39').Value = 'Ronald F. VanSteeland was the Vice President for Finance and Administration at GVSU serving between 1966 and 2001.  He helped build the financial backbone of GVSU and make it one of the “Top 100 Best Buys” in the US.  He was a board member of the West Michigan Environmental Action Council between 2006 and 2008.  He also met his wife Peggy while working at GVSU – she was secretary to President Don Lubbers.'
$ws.Range('G39').VerticalAlignment = -4108

$ws.Range('G40').Value = 'Ella Koeze Weed was born in in Grand Rapids, MI and graduated from South High School, Davenport Business College.  She entered politics in 1940 and is known as a frontrunner for all women to follow in politics.  She was elected to Michigan’s Constitutional Convention in 1961 and was appointed Chairman of National Alliance of Businessmen in 1971.  She also served on GVSU Board of Control between 1967 and 1981.  She married Albert Koeze and helped manage his successful peanut butter factory.  '
$ws.Range('G40').VerticalAlignment = -4108

$ws.Range('G41').Value = 'James H. Zumberge was the president of Grand Valley State Colleges between 1963 and 1969.  After that he became president of Southern Methodist University and fourth president of University of California.  As an educator he was a geology professor, even leading expedition to Antarctica.  Cape Zumberge and Zumberge Coast in Antarctica are named for him.'
$ws.Range('G41').VerticalAlignment = -4108

$ws.Range('G42').Value = 'Richard M. DeVos was the president and chairman of the Grand valley University Foundation.  He is the co-founder of Amway Corporation.  He is both an entrepreneur and a philanthropist.  He has a love for sports, owning the Orlando Magic, and NBA franchise in Florida, naming the arena for the company he started.  He was instrumental in leading the revitalization of downtown Grand Rapids when he opened the Amway Grand plaza Hotel.  He is one of Grand Valleys’ most charitable donors leading nearly every campaign since the early 1990’s.  '
$ws.Range('G42').VerticalAlignment = -4108

$ws.Range('G43').Value = 'Peter C. Cook created Mazda Great Lakes, the first and only independent direct distributorship in the U.S. until Mazda bought it back.  He is a renowned member of the prestigious Grand Valley Hall of Fame.  He also enjoyed playing the trumpet, as well as he had a love for Dixie bands.  He believed strongly in giving back to his community stating that, “There are two necessary things in giving: one is to be willing, the other is to be able.  Able is secondary because many people who are not willing to do become able but still never give.  People who have little yet give what they can become the big givers because they are able.  That’s because they had always been willing.”'
$ws.Range('G43').VerticalAlignment = -4108

$ws.Range('G44').Value = 'Peter C. Cook created Mazda Great Lakes, the first and only independent direct distributorship in the U.S. until Mazda bought it back.  He is a renowned member of the prestigious Grand Valley Hall of Fame.  He also enjoyed playing the trumpet, as well as he had a love for Dixie bands.  He believed strongly in giving back to his community stating that, “There are two necessary things in giving: one is to be willing, the other is to be able.  Able is secondary because many people who are not willing to do become able but still never give.  People who have little yet give what they can become the big givers because they are able.  That’s because they had always been willing.”'
$ws.Range('G44').VerticalAlignment = -4108

$ws.Range('G45').Value = 'Marvin Dewitt started Bil Mar Foods when he was 18 years old and bought 17 turkeys for $60.  Bill Mar was one of the top three turkey producers in the country.  He served on Zeeland School Board, Ottawa County Road Commission and the board of Northwestern College in Iowa.  He also served on the Grand Valley University Foundation’s Advisory Cabinet between 2001 and 2009.  '
$ws.Range('G45').VerticalAlignment = -4108

$ws.Range('G46').Value = 'Paul B. Henry was elected to the Michigan State House of Representatives between 1979 and 1982 and served on the Michigan State Senate between 1983 and 1985.  He was elected to U.S. House of Representatives in 1984, served until his death in 1993.  He was also a Peace Corps volunteer in Liberia and Ethiopia between 1963 and 1965.'
$ws.Range('G46').VerticalAlignment = -4108

# Update the sheet view: scroll so column F is the leftmost visible column, and select G46
$ws.Activate()
$ws.Range('G46').Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

